# Scheduled-runner update: refresh computed profit figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8749.25
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7472
$ws.Range("I62").Value = 6597
$ws.Range("J62").Value = 9222
$ws.Range("K62").Value = 6597
$ws.Range("L62").Value = 9222
$ws.Range("M62").Value = -5973
$ws.Range("N62").Value = -10470

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 7472
$ws.Range("I65").Value = 6597
$ws.Range("J65").Value = 9222
$ws.Range("K65").Value = 32985
$ws.Range("L65").Value = 46110
$ws.Range("M65").Value = -29865
$ws.Range("N65").Value = -52350

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3590
$ws.Range("I74").Value = 3590
$ws.Range("K74").Value = 3590
$ws.Range("M74").Value = -2654

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3590
$ws.Range("I77").Value = 3590
$ws.Range("K77").Value = 17950
$ws.Range("M77").Value = -13270

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6185.3145
$ws.Range("I116").Value = 5599.1333
$ws.Range("K116").Value = 5599.1333
$ws.Range("M116").Value = -2157.1333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8474.299999999999
$ws.Range("I138").Value = 6782.222
$ws.Range("K138").Value = 20346.666
$ws.Range("M138").Value = -15206.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 999
$ws.Range("I6").Value = 999
$ws.Range("K6").Value = 999
$ws.Range("M6").Value = -826

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1091.2
$ws.Range("I25").Value = 412.66666
$ws.Range("J25").Value = 2109
$ws.Range("K25").Value = 412.66666
$ws.Range("L25").Value = 2109
$ws.Range("M25").Value = -10.66665999999998
$ws.Range("N25").Value = -2913

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2499.5
$ws.Range("I74").Value = 2499.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2499.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1625.5
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2499.5
$ws.Range("I77").Value = 2499.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 12497.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -8129.5
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2063.3635
$ws.Range("I102").Value = 1188.5555
$ws.Range("K102").Value = 1188.5555
$ws.Range("M102").Value = 433.4445000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 925
$ws.Range("I132").Value = 925
$ws.Range("K132").Value = 2775
$ws.Range("M132").Value = -245

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3024.7778
$ws.Range("I86").Value = 2452.6667
$ws.Range("K86").Value = 2452.6667
$ws.Range("M86").Value = -1329.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3024.7778
$ws.Range("I89").Value = 2452.6667
$ws.Range("K89").Value = 12263.3335
$ws.Range("M89").Value = -6647.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2724.2104
$ws.Range("I99").Value = 2166.1538
$ws.Range("J99").Value = 3933.3333
$ws.Range("K99").Value = 2166.1538
$ws.Range("L99").Value = 3933.3333
$ws.Range("M99").Value = -668.1538
$ws.Range("N99").Value = -6929.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 2500
$ws.Range("J13").Value = 2500
$ws.Range("L13").Value = 2500
$ws.Range("N13").Value = -2778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 165.5
$ws.Range("I2").Value = 186.4
$ws.Range("J2").Value = 130.66667
$ws.Range("K2").Value = 1118.4
$ws.Range("L2").Value = 784.0000200000001
$ws.Range("M2").Value = -1005.4
$ws.Range("N2").Value = -1010.00002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 126.454544
$ws.Range("I38").Value = 50.125
$ws.Range("J38").Value = 330
$ws.Range("K38").Value = 150.375
$ws.Range("L38").Value = 990
$ws.Range("M38").Value = 196.625
$ws.Range("N38").Value = -1684

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1499.875
$ws.Range("I68").Value = 1759.8
$ws.Range("J68").Value = 1066.6666
$ws.Range("K68").Value = 5279.4
$ws.Range("L68").Value = 3199.9998
$ws.Range("M68").Value = -4468.4
$ws.Range("N68").Value = -4821.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1499.875
$ws.Range("I71").Value = 1759.8
$ws.Range("J71").Value = 1066.6666
$ws.Range("K71").Value = 15838.2
$ws.Range("L71").Value = 9599.999400000001
$ws.Range("M71").Value = -11782.2
$ws.Range("N71").Value = -17711.9994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 1499
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4570
$ws.Range("I129").Value = 3500
$ws.Range("J129").Value = 4998
$ws.Range("K129").Value = 10500
$ws.Range("L129").Value = 14994
$ws.Range("M129").Value = -5500
$ws.Range("N129").Value = -24994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3295
$ws.Range("I132").Value = 1997
$ws.Range("K132").Value = 17973
$ws.Range("M132").Value = -15443

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 751265.9399999999
$ws.Range("I11").Value = 751265.9399999999
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 751265.9399999999
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -751126.9399999999
$ws.Range("N11").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 49827.57
$ws.Range("I122").Value = 57638.723
$ws.Range("J122").Value = 2960.6667
$ws.Range("K122").Value = 172916.169
$ws.Range("L122").Value = 8882.000100000001
$ws.Range("M122").Value = -170466.169
$ws.Range("N122").Value = -13782.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7999
$ws.Range("I136").Value = 2997
$ws.Range("K136").Value = 8991
$ws.Range("M136").Value = -6441

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 10000
$ws.Range("I14").Value = 10000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -9832
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2846.8823
$ws.Range("I132").Value = 2796.0667
$ws.Range("K132").Value = 8388.2001
$ws.Range("M132").Value = -5858.2001
